# Natmi following Dr Hou advice
# The LR-pair table (Fgf2 -> Gpc4) is recomputed: an "ECs" sending cluster
# is added (previously only FAPs/sCs sent), so the 2 sender x 3 target
# (6-row) table becomes a 3 sender x 3 target (9-row) table, and every
# edge statistic is recalculated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$data = @(
    @("ECs",  "Fgf2", "Gpc4", "ECs",  2, 0.6666666666666666, 0.6462393333333333, 1.938718,          0.03461850536298827, 0.03461850536298827, 3, 1, 8.767882,          26.303646,  0.2073466531282026, 0.2073466531282026, 5.666150218425333,  50.995351965828,    0.007178031223316352, 0.00717803122331635),
    @("ECs",  "Fgf2", "Gpc4", "FAPs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718,          0.03461850536298827, 0.03461850536298827, 3, 1, 20.88867166666667, 62.666015,  0.4939843121037951, 0.4939843121037952, 13.49908125208555,  121.49173126877,    0.0171009985577973,   0.0171009985577973),
    @("ECs",  "Fgf2", "Gpc4", "sCs",  2, 0.6666666666666666, 0.6462393333333333, 1.938718,          0.03461850536298827, 0.03461850536298827, 3, 1, 12.62954966666667, 37.888649,  0.2986690347680022, 0.2986690347680022, 8.16171175688689,   73.455405811982,    0.01033947558187461,  0.01033947558187461),
    @("FAPs", "Fgf2", "Gpc4", "ECs",  3, 1,                   15.322826,           45.968478,        0.8208310864042159, 0.8208310864042158, 3, 1, 8.767882,          26.303646,  0.2073466531282026, 0.2073466531282026, 134.348730274532,   1209.138572470788,  0.1701965785495007,   0.1701965785495007),
    @("FAPs", "Fgf2", "Gpc4", "FAPs", 3, 1,                   15.322826,           45.968478,        0.8208310864042159, 0.8208310864042158, 3, 1, 20.88867166666667, 62.666015,  0.4939843121037951, 0.4939843121037952, 320.0734813194634,  2880.66133187517,   0.4054776795707974,   0.4054776795707974),
    @("FAPs", "Fgf2", "Gpc4", "sCs",  3, 1,                   15.322826,           45.968478,        0.8208310864042159, 0.8208310864042158, 3, 1, 12.62954966666667, 37.888649,  0.2986690347680022, 0.2986690347680022, 193.5203920006913,  1741.683528006222,  0.2451568282839177,   0.2451568282839177),
    @("sCs",  "Fgf2", "Gpc4", "ECs",  3, 1,                   2.698388,            8.095164,         0.1445504082327959, 0.1445504082327959, 3, 1, 8.767882,          26.303646,  0.2073466531282026, 0.2073466531282026, 23.659147574216,    212.932328167944,   0.02997204335538561,  0.02997204335538561),
    @("sCs",  "Fgf2", "Gpc4", "FAPs", 3, 1,                   2.698388,            8.095164,         0.1445504082327959, 0.1445504082327959, 3, 1, 20.88867166666667, 62.666015,  0.4939843121037951, 0.4939843121037952, 56.36574096127333,  507.29166865146,    0.07140563397520044,  0.07140563397520044),
    @("sCs",  "Fgf2", "Gpc4", "sCs",  3, 1,                   2.698388,            8.095164,         0.1445504082327959, 0.1445504082327959, 3, 1, 12.62954966666667, 37.888649,  0.2986690347680022, 0.2986690347680022, 34.07942526593734,  306.7148273934361,  0.04317273090220982,  0.04317273090220983)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
